# Updated cryptos list - applies price (D) and volume/1h change (E) updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.792.11"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "2.413.18"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.49"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.02"
$ws.Range("E6").Value = "  +4.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("D9").Value = "2.428.00"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  +3.80%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("E12").Value = "  +3.93%  "
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").Value = "2.835.75"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "61.685.39"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "2.420.75"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.74"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.29"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.04"
$ws.Range("E23").Value = "  +12.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.06"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.19"
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "621.21"
$ws.Range("E27").Value = "  +12.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.39"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").Value = "0.0₃0961"
$ws.Range("E29").Value = "  +4.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.06"
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("E32").Value = "  +8.83%  "
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("E35").Value = "  +5.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.371"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.35"
$ws.Range("E40").Value = "  +5.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.44"
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("E42").Value = "  +11.00%  "
$ws.Range("E43").Value = "  +4.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.32"
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "0.0₆0284"
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.36"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.00"
$ws.Range("E49").Value = "  +4.73%  "
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0510"
$ws.Range("E51").Value = "  +2.25%  "
